# Applies the vocabulary.xlsx edits described by the commit:
# "new .ttl from Google sheet has been generated"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Namespace URI changed from test2 to test3 ---
$ws.Range("B1").Value = "http://purl.org/test3/variables/"
$ws.Range("C3").Value = "http://purl.org/test3/variables/"

# --- Vocabulary title / description simplified to "test" ---
$ws.Range("B10").Value = "test"
$ws.Range("B11").Value = "test"

# --- Row 13: now holds dct:creator / Minka (was dct:rights / license text) ---
$ws.Range("A13").Value = "dct:creator"
$ws.Range("B13").Value = "Minka"
$ws.Range("C13").Value = ""

# --- Row 14: now holds dct:rights / license text (was pav:version) ---
$ws.Range("A14").Value = "dct:rights"
$ws.Range("C14").Value = "License under which the vocabulary is provided"

# --- Row 15: now holds pav:version (was pav:createdOn) ---
$ws.Range("A15").Value = "pav:version"
$ws.Range("C15").Value = "Vocabulary version"

# --- Row 16: now holds pav:createdOn (was pav:lastUpdatedOn) ---
$ws.Range("A16").Value = "pav:createdOn"
$ws.Range("C16").Value = "Date when vocabulary was initially created (follow https://en.wikipedia.org/wiki/ISO_8601)"

# --- Row 17: now holds pav:lastUpdatedOn (was "Definition of terms" heading) ---
$ws.Range("A17").Value = "pav:lastUpdatedOn"
$ws.Range("C17").Value = "Date of the last vocabulary update"

# --- Row 18: now holds the "Definition of terms" heading (was the Identifier header row) ---
$ws.Range("A18").Value = "Definition of terms (optionally properties)"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""
$ws.Range("J18").Value = ""
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""

# --- Row 19: now holds the column-header definitions (was "vars:Test" sample row) ---
$ws.Range("A19").Value = "Identifier"
$ws.Range("B19").Value = "skos:prefLabel@en"
$ws.Range("C19").Value = 'qudt:unit(separator=",")'
$ws.Range("D19").Value = 'skos:altLabel(separator=";")'
$ws.Range("E19").Value = "skos:definition@en"
$ws.Range("F19").Value = 'dct:source(separator=",")'
$ws.Range("G19").Value = 'skos:broader(separator=",")'
$ws.Range("H19").Value = 'skos:exactMatch(separator=",")'
$ws.Range("I19").Value = 'skos:closeMatch(separator=",")'
$ws.Range("J19").Value = "skos:editorialNote@en"
$ws.Range("K19").Value = 'dct:creator(separator=",")'
$ws.Range("L19").Value = 'dct:contributor(separator=",")'

# --- Row 20: renamed term vars:test / test, definition & broader cleared ---
$ws.Range("A20").Value = "vars:test"
$ws.Range("B20").Value = "test"
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""

# --- Row 21: renamed term vars:computerscientist with a new definition ---
$ws.Range("A21").Value = "vars:computerscientist"
$ws.Range("B21").Value = "computerscientist"
$ws.Range("E21").Value = "a person that knows stuff about computers"

# --- New row 90 appended at the bottom, following the same "vars:" placeholder pattern as rows 22-89 ---
$ws.Range("A90").Value = "vars:"
